# Add a new row (row 11) with 2021 data to Sheet1, matching the style of the
# previous year row (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing data row (row 10) down to the new
# row (row 11) so the year label keeps the same bold/bordered/centered style.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

# Year label for the new row
$ws.Range("A11").Value = "2021年"

# Data values for each industry column (row 11 = 2021年)
$ws.Range("B11").Value = 10380.07
$ws.Range("C11").Value = 2440.45
$ws.Range("D11").Value = 2002.71

# E11 stays blank (matches the other years in this column, which are empty)
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats -> creates the blank cell without a value

$ws.Range("F11").Value = 2741.68
$ws.Range("G11").Value = 25343.09
$ws.Range("H11").Value = 1730.98
$ws.Range("I11").Value = 7227.07
$ws.Range("J11").Value = 1033.5
$ws.Range("K11").Value = 565082.11
$ws.Range("L11").Value = 195.59
$ws.Range("M11").Value = 753.72
$ws.Range("N11").Value = 2621.95
$ws.Range("O11").Value = 351.15
$ws.Range("P11").Value = 16617.01
$ws.Range("Q11").Value = 3710.26
$ws.Range("R11").Value = 277.1
$ws.Range("S11").Value = 1906.46
$ws.Range("T11").Value = 19917.64
$ws.Range("U11").Value = 37419.49
$ws.Range("V11").Value = 11621.32
$ws.Range("W11").Value = 52182.19
$ws.Range("X11").Value = 7934.76
$ws.Range("Y11").Value = 171657.4
$ws.Range("Z11").Value = 9041.83
$ws.Range("AA11").Value = 80.06999999999999
$ws.Range("AB11").Value = 20517.9
$ws.Range("AC11").Value = 19468.78
$ws.Range("AD11").Value = 1037.12
$ws.Range("AE11").Value = 414.52
$ws.Range("AF11").Value = 26065.27
$ws.Range("AG11").Value = 12141.12
$ws.Range("AH11").Value = 1493.36
$ws.Range("AI11").Value = 8277.75
$ws.Range("AJ11").Value = 1696.04
$ws.Range("AK11").Value = 5472.36
$ws.Range("AL11").Value = 21343.16
$ws.Range("AM11").Value = 15455.3
$ws.Range("AN11").Value = 3121.78
$ws.Range("AO11").Value = 1911.55
$ws.Range("AP11").Value = 28881.84
$ws.Range("AQ11").Value = 8589.84
